# "Fixed weird visual bug on presentation"
#
# 1) Slide 12 title placeholder was too narrow (text was clipping/wrapping
#    oddly) -> widen it from 5478600 EMU to 5621700 EMU (height/position
#    unchanged).
# 2) The deck's theme (theme1.xml, used by the slide master / all slides)
#    was carrying the wrong ("Simple Light") accent/background palette;
#    the intended palette is the "Default" one. Re-point every theme
#    color slot to the correct swatch so the slides render with the
#    right colors.

$p = $ppt.ActivePresentation

# --- 1) Widen the title box on slide 12 -----------------------------------
$s  = $p.Slides.Item(12)
$sh = $s.Shapes.Item(1)
# Target width is 5621700 EMU (1 pt = 12700 EMU => 442.653543307... pt).
# Nudged a hair above the exact quotient so the internal EMU conversion
# lands exactly on 5621700 rather than rounding down.
$sh.Width = 442.65357971191406

# --- 2) Swap the theme's color scheme onto the correct palette ------------
$master = $p.SlideMaster
$cs = $master.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0          # dk1      000000
$cs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Item(3).RGB  = 5800213    # dk2      158158
$cs.Item(4).RGB  = 15987699   # lt2      F3F3F3
$cs.Item(5).RGB  = 13077765   # accent1  058DC7
$cs.Item(6).RGB  = 3322960    # accent2  50B432
$cs.Item(7).RGB  = 1791725    # accent3  ED561B
$cs.Item(8).RGB  = 61421      # accent4  EDEF00
$cs.Item(9).RGB  = 15059748   # accent5  24CBE5
$cs.Item(10).RGB = 7529828    # accent6  64E572
$cs.Item(11).RGB = 13369378   # hlink    2200CC
$cs.Item(12).RGB = 9116245    # folHlink 551A8B
